# Fruta / hortaliza, semanal
# Insert two new rows of data at the top of the data block (rows 1176-1177),
# pushing the existing rows 1176-1274 down to 1178-1276.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 1176 (existing row 1176 and everything below
# shifts down by two rows, becoming rows 1178-1276).
$ws.Range("A1176:A1177").EntireRow.Insert()

# --- New row 1176 ---
$ws.Cells.Item(1176, 1).Value = 10
$ws.Cells.Item(1176, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1176, 3).Value = "La Araucanía"
$ws.Cells.Item(1176, 4).Value = 45013
$ws.Cells.Item(1176, 5).Value = 9
$ws.Cells.Item(1176, 6).Value = "Fruta"
$ws.Cells.Item(1176, 7).Value = 100109
$ws.Cells.Item(1176, 8).Value = "Uva"
$ws.Cells.Item(1176, 9).Value = 100109001
$ws.Cells.Item(1176, 10).Value = "Uva"
$ws.Cells.Item(1176, 11).Value = "Red Globe"
$ws.Cells.Item(1176, 12).Value = "Primera"
$ws.Cells.Item(1176, 13).Value = 125
$ws.Cells.Item(1176, 14).Value = 12000
$ws.Cells.Item(1176, 15).Value = 12000
$ws.Cells.Item(1176, 16).Value = 12000
$ws.Cells.Item(1176, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(1176, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1176, 19).Value = 667
$ws.Cells.Item(1176, 20).Value = 18

# --- New row 1177 ---
$ws.Cells.Item(1177, 1).Value = 10
$ws.Cells.Item(1177, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1177, 3).Value = "La Araucanía"
$ws.Cells.Item(1177, 4).Value = 45013
$ws.Cells.Item(1177, 5).Value = 9
$ws.Cells.Item(1177, 6).Value = "Fruta"
$ws.Cells.Item(1177, 7).Value = 100109
$ws.Cells.Item(1177, 8).Value = "Uva"
$ws.Cells.Item(1177, 9).Value = 100109001
$ws.Cells.Item(1177, 10).Value = "Uva"
$ws.Cells.Item(1177, 11).Value = "Thompson seedless"
$ws.Cells.Item(1177, 12).Value = "Primera"
$ws.Cells.Item(1177, 13).Value = 215
$ws.Cells.Item(1177, 14).Value = 12000
$ws.Cells.Item(1177, 15).Value = 12000
$ws.Cells.Item(1177, 16).Value = 12000
$ws.Cells.Item(1177, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(1177, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1177, 19).Value = 667
$ws.Cells.Item(1177, 20).Value = 18
